$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.382.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.883.01"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7127"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.39"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08042"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3125"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.31"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08331"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.906.27"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.247"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7200"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.336"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +5.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008516"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.393.84"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.59"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.142.73"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.855"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1586"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.17%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.066"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.420"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.348"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05368"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.949"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7499"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.699"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01890"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.288.75"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.741"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.617"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9237"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "111.92"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.33"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.047.29"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.87%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.807"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5220"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.516"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4393"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.133"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.99%  "
